$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 74

$ws.Cells.Item($newRow, 1).Value = "December"
$ws.Cells.Item($newRow, 2).Value = 17
$ws.Cells.Item($newRow, 3).Value = "19:52:08"
$ws.Cells.Item($newRow, 4).Value = 1.13
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 70
$ws.Cells.Item($newRow, 7).Value = 1031
$ws.Cells.Item($newRow, 8).Value = 2.24
